$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.102.35'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '1.652.88'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = '218.59'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').Value = '0.5299'
$ws.Range('E6').Value = '  +1.28%  '
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').Value = '0.2609'
$ws.Range('E8').Value = '  -2.19%  '
$ws.Range('D9').Value = '0.06327'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').Value = '20.40'
$ws.Range('E10').Value = '  -2.98%  '
$ws.Range('D11').Value = '0.07746'
$ws.Range('E11').Value = '  +0.46%  '
$ws.Range('D12').Value = '4.487'
$ws.Range('E12').Value = '  +1.35%  '
$ws.Range('D13').Value = '1.649.94'
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('D14').Value = '0.5462'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').Value = '0.0₅8115'
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('D16').Value = '65.25'
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('D17').Value = '26.118.57'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').Value = '1.003'
$ws.Range('D19').Value = '4.545'
$ws.Range('E19').Value = '  -2.52%  '
$ws.Range('D20').Value = '194.01'
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('E22').Value = '  -1.50%  '
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('E24').Value = '  +0.75%  '
$ws.Range('D25').Value = '0.1241'
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('D26').Value = '7.266'
$ws.Range('E26').Value = '  +0.47%  '
$ws.Range('D27').Value = '16.16'
$ws.Range('E27').Value = '  +0.18%  '
$ws.Range('E28').Value = '  +1.27%  '
$ws.Range('D29').Value = '0.05938'
$ws.Range('E29').Value = '  -0.96%  '
$ws.Range('D30').Value = '1.279'
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('D31').Value = '3.506'
$ws.Range('E31').Value = '  -5.00%  '
$ws.Range('D32').Value = '3.238'
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('D33').Value = '1.549'
$ws.Range('E33').Value = '  -5.16%  '
$ws.Range('D34').Value = '2.412'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('D35').Value = '0.9446'
$ws.Range('E35').Value = '  -3.57%  '
$ws.Range('D36').Value = '2.756'
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').Value = '0.5637'
$ws.Range('E37').Value = '  -3.83%  '
$ws.Range('D38').Value = '0.01607'
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').Value = '5.857'
$ws.Range('E39').Value = '  -1.48%  '
$ws.Range('D40').Value = '0.8465'
$ws.Range('E40').Value = '  -1.64%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').Value = '100.84'
$ws.Range('E42').Value = '  +1.19%  '
$ws.Range('D43').Value = '1.008.16'
$ws.Range('E43').Value = '  -2.47%  '
$ws.Range('D44').Value = '1.799.37'
$ws.Range('E44').Value = '  -0.23%  '
$ws.Range('E46').Value = '  -0.95%  '
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('D48').Value = '0.4288'
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('D49').Value = '0.05152'
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('D50').Value = '1.470'
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('D51').Value = '7.769'
$ws.Range('E51').Value = '  -3.79%  '
